# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.080.72"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.73%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.858.39"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +3.24%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.27"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +4.29%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +2.16%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.55"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +9.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.328"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.45%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.90%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.37%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +3.17%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.39"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.856.15"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.68%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.679"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.61%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "35.052.97"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.48%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.02%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.66%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "241.02"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +3.84%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.28"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "171.47"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +30.94%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.13%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.69"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +3.51%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.39%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0559"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +4.01%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.74%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.01"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.03"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +14.41%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +23.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.30"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +5.54%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.780"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +13.31%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +13.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "91.94"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0203"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +6.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.352.18"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.92%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.90"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +4.85%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +6.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.70"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +56.08%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0547"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +6.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.35"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +2.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.043.96"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.18%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.43"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +17.93%  "
